$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.64"
$ws.Range("E2").Value = "'-4.17%"
$ws.Range("D3").Value = "'30.91"
$ws.Range("E3").Value = "'-3.87%"
$ws.Range("D4").Value = "'4.947"
$ws.Range("E4").Value = "'0.12%"
$ws.Range("D5").Value = "'0.07159"
$ws.Range("E5").Value = "'-8.64%"
$ws.Range("D6").Value = "'1.800"
$ws.Range("E6").Value = "'-11.70%"
$ws.Range("D7").Value = "'7.664"
$ws.Range("E7").Value = "'-2.27%"
$ws.Range("D8").Value = "'3.732"
$ws.Range("E8").Value = "'-2.97%"
$ws.Range("D9").Value = "'0.8949"
$ws.Range("E9").Value = "'-2.98%"
$ws.Range("D10").Value = "'0.1652"
$ws.Range("E10").Value = "'-6.19%"
$ws.Range("D11").Value = "'0.07669"
$ws.Range("E11").Value = "'-2.97%"
$ws.Range("D12").Value = "'0.08090"
$ws.Range("E12").Value = "'-5.81%"
$ws.Range("D13").Value = "'0.03028"
$ws.Range("E13").Value = "'-4.19%"
$ws.Range("D14").Value = "'0.1003"
$ws.Range("E14").Value = "'-0.16%"
$ws.Range("D15").Value = "'0.001497"
$ws.Range("E15").Value = "'-0.76%"
$ws.Range("D16").Value = "'0.005759"
$ws.Range("E16").Value = "'0.10%"
$ws.Range("D17").Value = "'3.477"
$ws.Range("E17").Value = "'0.35%"
$ws.Range("D18").Value = "'2.083"
$ws.Range("E18").Value = "'-3.41%"
$ws.Range("E19").Value = "'0.04%"
$ws.Range("D20").Value = "'0.1272"
$ws.Range("E20").Value = "'-3.42%"
$ws.Range("D21").Value = "'4.036"
$ws.Range("E21").Value = "'-5.87%"
$ws.Range("E22").Value = "'0.36%"
$ws.Range("D23").Value = "'0.04509"
$ws.Range("E23").Value = "'-1.39%"
$ws.Range("D24").Value = "'0.001212"
$ws.Range("E24").Value = "'-0.93%"
$ws.Range("D25").Value = "'0.004006"
$ws.Range("E25").Value = "'-10.00%"
$ws.Range("D26").Value = "'0.0001249"
$ws.Range("E26").Value = "'-0.07%"
$ws.Range("D39").Value = "'0.01604"
$ws.Range("E39").Value = "'-8.08%"
$ws.Range("D40").Value = "'0.04385"
$ws.Range("E40").Value = "'-8.39%"
$ws.Range("D41").Value = "'0.007345"
$ws.Range("E41").Value = "'-2.76%"
$ws.Range("D42").Value = "'0.1305"
$ws.Range("E42").Value = "'-4.51%"
$ws.Range("D43").Value = "'0.002004"
$ws.Range("E43").Value = "'-15.08%"
$ws.Range("D44").Value = "'0.009226"
$ws.Range("E44").Value = "'-13.16%"
$ws.Range("D45").Value = "'0.00005944"
$ws.Range("E45").Value = "'-5.88%"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("E46").Value = "'-0.14%"
$ws.Range("D47").Value = "'2.247"
$ws.Range("E47").Value = "'173.92%"
$ws.Range("D48").Value = "'0.002996"
$ws.Range("E48").Value = "'-3.41%"
$ws.Range("D49").Value = "'0.00002097"
$ws.Range("E49").Value = "'-0.14%"
$ws.Range("D50").Value = "'0.0001997"
$ws.Range("E50").Value = "'-0.14%"
